# Scheduled market-data refresh: update cached Universalis price columns
# (H..N) on the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1510.4912
$ws.Range("I15").Value = 1510.4912
$ws.Range("K15").Value = 4531.473599999999
$ws.Range("M15").Value = -4362.473599999999

$ws.Range("H28").Value = 761.17645
$ws.Range("I28").Value = 389
$ws.Range("K28").Value = 389
$ws.Range("M28").Value = 96

$ws.Range("H74").Value = 5654.8184
$ws.Range("I74").Value = 4717.1665
$ws.Range("J74").Value = 6780
$ws.Range("K74").Value = 4717.1665
$ws.Range("L74").Value = 6780
$ws.Range("M74").Value = -3781.1665
$ws.Range("N74").Value = -8652

$ws.Range("H77").Value = 5654.8184
$ws.Range("I77").Value = 4717.1665
$ws.Range("J77").Value = 6780
$ws.Range("K77").Value = 23585.8325
$ws.Range("L77").Value = 33900
$ws.Range("M77").Value = -18905.8325
$ws.Range("N77").Value = -43260

$ws.Range("H103").Value = 773.3333
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 773.3333
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2319.9999
$ws.Range("M103").Value = ""
$ws.Range("N103").Value = -3491.9999

$ws.Range("H129").Value = 939.2895
$ws.Range("I129").Value = 726.8
$ws.Range("J129").Value = 1015.1786
$ws.Range("K129").Value = 2180.4
$ws.Range("L129").Value = 3045.5358
$ws.Range("M129").Value = 2819.6
$ws.Range("N129").Value = -13045.5358

$ws.Range("H132").Value = 2930.9062
$ws.Range("I132").Value = 3052.9666
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 9158.899800000001
$ws.Range("L132").Value = 3300
$ws.Range("M132").Value = -6628.899800000001
$ws.Range("N132").Value = -8360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2290.0952
$ws.Range("J61").Value = 2864.2856
$ws.Range("L61").Value = 2864.2856
$ws.Range("N61").Value = -3288.2856

$ws.Range("H122").Value = 1892.3572
$ws.Range("I122").Value = 1678.8422
$ws.Range("J122").Value = 2343.111
$ws.Range("K122").Value = 5036.5266
$ws.Range("L122").Value = 7029.333
$ws.Range("M122").Value = -2586.5266
$ws.Range("N122").Value = -11929.333

$ws.Range("H136").Value = 2290.0952
$ws.Range("J136").Value = 2864.2856
$ws.Range("L136").Value = 8592.856800000001
$ws.Range("N136").Value = -13692.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 135371.8
$ws.Range("I105").Value = 126702.375
$ws.Range("K105").Value = 126702.375
$ws.Range("M105").Value = -124955.375

$ws.Range("H108").Value = 36000
$ws.Range("J108").Value = 36000
$ws.Range("L108").Value = 36000
$ws.Range("N108").Value = -43680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1798.3334
$ws.Range("I31").Value = 1130.5
$ws.Range("K31").Value = 1130.5
$ws.Range("M31").Value = -835.5

$ws.Range("H34").Value = 1798.3334
$ws.Range("I34").Value = 1130.5
$ws.Range("K34").Value = 1130.5
$ws.Range("M34").Value = -928.5

$ws.Range("H62").Value = 5053351.5
$ws.Range("I62").Value = 18520686
$ws.Range("J62").Value = 3100.75
$ws.Range("K62").Value = 18520686
$ws.Range("L62").Value = 3100.75
$ws.Range("M62").Value = -18520062
$ws.Range("N62").Value = -4348.75

$ws.Range("H65").Value = 5053351.5
$ws.Range("I65").Value = 18520686
$ws.Range("J65").Value = 3100.75
$ws.Range("K65").Value = 92603430
$ws.Range("L65").Value = 15503.75
$ws.Range("M65").Value = -92600310
$ws.Range("N65").Value = -21743.75

$ws.Range("H132").Value = 2610.3157
$ws.Range("I132").Value = 2677.0833
$ws.Range("J132").Value = 2495.8572
$ws.Range("K132").Value = 8031.249899999999
$ws.Range("L132").Value = 7487.571599999999
$ws.Range("M132").Value = -5501.249899999999
$ws.Range("N132").Value = -12547.5716

$ws.Range("H134").Value = 1321.5
$ws.Range("I134").Value = 1187.3684
$ws.Range("J134").Value = 2171
$ws.Range("K134").Value = 3562.1052
$ws.Range("L134").Value = 6513
$ws.Range("M134").Value = -1027.1052
$ws.Range("N134").Value = -11583

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 400
$ws.Range("I63").Value = 400
$ws.Range("K63").Value = 1200
$ws.Range("M63").Value = -451

$ws.Range("H64").Value = 2141.2856
$ws.Range("I64").Value = 996.3333
$ws.Range("K64").Value = 2988.9999
$ws.Range("M64").Value = -2718.9999

$ws.Range("H66").Value = 400
$ws.Range("I66").Value = 400
$ws.Range("K66").Value = 3600
$ws.Range("M66").Value = 144

$ws.Range("H67").Value = 2141.2856
$ws.Range("I67").Value = 996.3333
$ws.Range("K67").Value = 2988.9999
$ws.Range("M67").Value = -2052.9999

$ws.Range("H107").Value = 248891.84
$ws.Range("I107").Value = 577.8095
$ws.Range("J107").Value = 397880.25
$ws.Range("K107").Value = 1733.4285
$ws.Range("L107").Value = 1193640.75
$ws.Range("M107").Value = 186.5715
$ws.Range("N107").Value = -1197480.75

$ws.Range("H131").Value = 882.53625
$ws.Range("J131").Value = 930.6229
$ws.Range("L131").Value = 2791.8687
$ws.Range("N131").Value = -12871.8687

$ws.Range("H136").Value = 5886.522
$ws.Range("I136").Value = 439.0909
$ws.Range("K136").Value = 1317.2727
$ws.Range("M136").Value = 3782.7273

$ws.Range("H138").Value = 1915.1875
$ws.Range("J138").Value = 3004.125
$ws.Range("L138").Value = 9012.375
$ws.Range("N138").Value = -19292.375

$ws.Range("H139").Value = 1923.138
$ws.Range("I139").Value = 669.4375
$ws.Range("J139").Value = 3466.1538
$ws.Range("K139").Value = 2008.3125
$ws.Range("L139").Value = 10398.4614
$ws.Range("M139").Value = 3131.6875
$ws.Range("N139").Value = -20678.4614

$ws.Range("H140").Value = 1738.3334
$ws.Range("I140").Value = 1441.8182
$ws.Range("K140").Value = 4325.4546
$ws.Range("M140").Value = 854.5454

$ws.Range("H141").Value = 2302.652
$ws.Range("I141").Value = 1275.6111
$ws.Range("K141").Value = 3826.8333
$ws.Range("M141").Value = 1353.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10723.333
$ws.Range("J52").Value = 10358.571
$ws.Range("L52").Value = 10358.571
$ws.Range("N52").Value = -10876.571

$ws.Range("H132").Value = 1977
$ws.Range("I132").Value = 1686.1666
$ws.Range("J132").Value = 2849.5
$ws.Range("K132").Value = 5058.4998
$ws.Range("L132").Value = 8548.5
$ws.Range("M132").Value = -2528.4998
$ws.Range("N132").Value = -13608.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 499.08334
$ws.Range("I55").Value = 159.86667
$ws.Range("K55").Value = 159.86667
$ws.Range("M55").Value = 13.13333

$ws.Range("H100").Value = 2978.2
$ws.Range("I100").Value = 2476.6
$ws.Range("J100").Value = 3479.8
$ws.Range("K100").Value = 2476.6
$ws.Range("L100").Value = 3479.8
$ws.Range("M100").Value = -1935.6
$ws.Range("N100").Value = -4561.8

$ws.Range("H132").Value = 3784.25
$ws.Range("I132").Value = 4014.92
$ws.Range("J132").Value = 3399.8
$ws.Range("K132").Value = 12044.76
$ws.Range("L132").Value = 10199.4
$ws.Range("M132").Value = -9514.76
$ws.Range("N132").Value = -15259.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 125002910
$ws.Range("I96").Value = 250002990
$ws.Range("J96").Value = 2844.75
$ws.Range("K96").Value = 250002990
$ws.Range("L96").Value = 2844.75
$ws.Range("M96").Value = -250001617
$ws.Range("N96").Value = -5590.75

$ws.Range("H107").Value = 47195.348
$ws.Range("I107").Value = 20439.6
$ws.Range("K107").Value = 61318.8
$ws.Range("M107").Value = -59398.8

$ws.Range("H122").Value = 976.5455
$ws.Range("I122").Value = 968.8571
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 2906.5713
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -456.5712999999996
$ws.Range("N122").Value = -7870

$ws.Range("H126").Value = 1911.5
$ws.Range("I126").Value = 1842.1111
$ws.Range("J126").Value = 2000.7142
$ws.Range("K126").Value = 5526.3333
$ws.Range("L126").Value = 6002.142599999999
$ws.Range("M126").Value = -3056.3333
$ws.Range("N126").Value = -10942.1426
